$wb = $excel.ActiveWorkbook

# wc_lang renamed its DfbaNetReaction/DfbaNetSpecies model classes to
# DfbaObjReaction/DfbaObjSpecies -- mirror that in the workbook's sheet
# names (renaming the sheet automatically keeps every defined name /
# autofilter reference that quotes the old sheet name in sync).
$wsNetReactions = $wb.Worksheets.Item("dFBA net reactions")
$wsNetReactions.Name = "dFBA objective reactions"

$wsNetSpecies = $wb.Worksheets.Item("dFBA net species")
$wsNetSpecies.Name = "dFBA objective species"

# The species sheet has a column header that spells out the (renamed)
# reaction sheet/concept -- update the text itself.
$wsNetSpecies.Range("C1").Value = "dFBA objective reaction"

# Leave the species sheet active/selected (as last edited), matching the
# workbook's persisted view state.
$wsNetSpecies.Activate()
$wsNetSpecies.Range("C2").Select()
